$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all edited cells keep their original text representation
# (values such as "1.001" or "0.07875" must not be reinterpreted as numbers)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.455.33'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.870.79'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.74'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7062'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3154'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07875'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08010'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.891.12'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.216'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.20'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7050'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.492'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008373'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.62%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.472.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '256.77'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.137.97'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.20'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.629'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.78%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1556'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.062'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.501'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.340'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.252'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05324'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7475'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.07%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.171'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.713'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01877'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.750'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8986'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.03'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.80'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.940'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -9.04%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.039.74'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.810'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.510'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06093'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.17%  '
